$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 924.35
$ws.Range("I15").Value = 924.35
$ws.Range("K15").Value = 2773.05
$ws.Range("M15").Value = -2604.05
$ws.Range("H42").Value = 62.125
$ws.Range("I42").Value = 12.25
$ws.Range("J42").Value = 112
$ws.Range("K42").Value = 36.75
$ws.Range("L42").Value = 336
$ws.Range("M42").Value = 193.25
$ws.Range("N42").Value = -796
$ws.Range("H80").Value = 1520.2778
$ws.Range("I80").Value = 1020.38464
$ws.Range("J80").Value = 2820
$ws.Range("K80").Value = 3061.15392
$ws.Range("L80").Value = 8460
$ws.Range("M80").Value = -2063.15392
$ws.Range("N80").Value = -10456
$ws.Range("H83").Value = 1520.2778
$ws.Range("I83").Value = 1020.38464
$ws.Range("J83").Value = 2820
$ws.Range("K83").Value = 9183.46176
$ws.Range("L83").Value = 25380
$ws.Range("M83").Value = -4191.46176
$ws.Range("N83").Value = -35364
$ws.Range("H100").Value = 2254
$ws.Range("I100").Value = 0
$ws.Range("K100").Value = 0
$ws.Range("M100").ClearContents()
$ws.Range("H113").Value = 1986.8
$ws.Range("I113").Value = 1764.1177
$ws.Range("J113").Value = 2460
$ws.Range("K113").Value = 1764.1177
$ws.Range("L113").Value = 2460
$ws.Range("M113").Value = 1489.8823
$ws.Range("N113").Value = -8968
$ws.Range("H116").Value = 1905.2941
$ws.Range("I116").Value = 1891.5385
$ws.Range("J116").Value = 1950
$ws.Range("K116").Value = 1891.5385
$ws.Range("L116").Value = 1950
$ws.Range("M116").Value = 1550.4615
$ws.Range("N116").Value = -8834
$ws.Range("H129").Value = 856.2708
$ws.Range("I129").Value = 516.7778
$ws.Range("J129").Value = 934.61536
$ws.Range("K129").Value = 1550.3334
$ws.Range("L129").Value = 2803.84608
$ws.Range("M129").Value = 3449.6666
$ws.Range("N129").Value = -12803.84608
$ws.Range("H137").Value = 25002824
$ws.Range("I137").Value = 866.3333
$ws.Range("K137").Value = 2598.9999
$ws.Range("M137").Value = -48.9998999999998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H33").Value = 14000
$ws.Range("I33").Value = 14000
$ws.Range("K33").Value = 14000
$ws.Range("M33").Value = -13671
$ws.Range("H37").Value = 7474.6
$ws.Range("I37").Value = 4000
$ws.Range("J37").Value = 8343.25
$ws.Range("K37").Value = 4000
$ws.Range("L37").Value = 8343.25
$ws.Range("M37").Value = -3727
$ws.Range("N37").Value = -8889.25
$ws.Range("H47").Value = 0
$ws.Range("J47").Value = 0
$ws.Range("L47").Value = 0
$ws.Range("N47").ClearContents()
$ws.Range("H88").Value = 5953.25
$ws.Range("I88").Value = 4403
$ws.Range("J88").Value = 7503.5
$ws.Range("K88").Value = 4403
$ws.Range("L88").Value = 7503.5
$ws.Range("M88").Value = -3997
$ws.Range("N88").Value = -8315.5
$ws.Range("H91").Value = 5953.25
$ws.Range("I91").Value = 4403
$ws.Range("J91").Value = 7503.5
$ws.Range("K91").Value = 4403
$ws.Range("L91").Value = 7503.5
$ws.Range("M91").Value = -2999
$ws.Range("N91").Value = -10311.5
$ws.Range("H127").Value = 41525.555
$ws.Range("I127").Value = 0
$ws.Range("J127").Value = 41525.555
$ws.Range("K127").Value = 0
$ws.Range("L127").Value = 41525.555
$ws.Range("M127").ClearContents()
$ws.Range("N127").Value = -51445.555
$ws.Range("H129").Value = 42694.75
$ws.Range("J129").Value = 42694.75
$ws.Range("L129").Value = 42694.75
$ws.Range("N129").Value = -52694.75
$ws.Range("H132").Value = 2175.5386
$ws.Range("I132").Value = 1980.4783
$ws.Range("J132").Value = 3671
$ws.Range("K132").Value = 5941.4349
$ws.Range("L132").Value = 11013
$ws.Range("M132").Value = -3411.4349
$ws.Range("N132").Value = -16073
$ws.Range("H137").Value = 76710
$ws.Range("I137").Value = 39900
$ws.Range("J137").Value = 79164
$ws.Range("K137").Value = 39900
$ws.Range("L137").Value = 79164
$ws.Range("M137").Value = -34800
$ws.Range("N137").Value = -89364

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 47785
$ws.Range("J53").Value = 47785
$ws.Range("L53").Value = 47785
$ws.Range("N53").Value = -48933
$ws.Range("H98").Value = 30542
$ws.Range("J98").Value = 30542
$ws.Range("L98").Value = 30542
$ws.Range("N98").Value = -36532
$ws.Range("H99").Value = 2498
$ws.Range("I99").Value = 0
$ws.Range("J99").Value = 2498
$ws.Range("K99").Value = 0
$ws.Range("L99").Value = 2498
$ws.Range("M99").ClearContents()
$ws.Range("N99").Value = -5494
$ws.Range("I107").Value = 700
$ws.Range("K107").Value = 700
$ws.Range("M107").Value = 1220
$ws.Range("H129").Value = 44191.145
$ws.Range("J129").Value = 44191.145
$ws.Range("L129").Value = 44191.145
$ws.Range("N129").Value = -54191.145
$ws.Range("H134").Value = 2849.4187
$ws.Range("I134").Value = 2554.1667
$ws.Range("K134").Value = 7662.500100000001
$ws.Range("M134").Value = -5127.500100000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10758448
$ws.Range("I31").Value = 8549
$ws.Range("J31").Value = 22225006
$ws.Range("K31").Value = 8549
$ws.Range("L31").Value = 22225006
$ws.Range("M31").Value = -8254
$ws.Range("N31").Value = -22225596
$ws.Range("H34").Value = 10758448
$ws.Range("I34").Value = 8549
$ws.Range("J34").Value = 22225006
$ws.Range("K34").Value = 8549
$ws.Range("L34").Value = 22225006
$ws.Range("M34").Value = -8347
$ws.Range("N34").Value = -22225410
$ws.Range("H35").Value = 2000
$ws.Range("I35").Value = 2000
$ws.Range("K35").Value = 2000
$ws.Range("M35").Value = -1706
$ws.Range("H39").Value = 8000
$ws.Range("I39").Value = 8000
$ws.Range("K39").Value = 8000
$ws.Range("M39").Value = -7609
$ws.Range("H47").Value = 29357
$ws.Range("J47").Value = 29357
$ws.Range("L47").Value = 29357
$ws.Range("N47").Value = -30489
$ws.Range("H49").Value = 8000
$ws.Range("I49").Value = 8000
$ws.Range("K49").Value = 8000
$ws.Range("M49").Value = -7818
$ws.Range("H132").Value = 33336040
$ws.Range("I132").Value = 41668940
$ws.Range("J132").Value = 4443.6665
$ws.Range("K132").Value = 125006820
$ws.Range("L132").Value = 13330.9995
$ws.Range("M132").Value = -125004290
$ws.Range("N132").Value = -18390.9995
$ws.Range("H134").Value = 723697
$ws.Range("I134").Value = 2174.0417
$ws.Range("J134").Value = 2647758.2
$ws.Range("K134").Value = 6522.125100000001
$ws.Range("L134").Value = 7943274.600000001
$ws.Range("M134").Value = -3987.125100000001
$ws.Range("N134").Value = -7948344.600000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1491.0968
$ws.Range("I5").Value = 288.35
$ws.Range("J5").Value = 3677.9092
$ws.Range("K5").Value = 865.0500000000001
$ws.Range("L5").Value = 11033.7276
$ws.Range("M5").Value = -753.0500000000001
$ws.Range("N5").Value = -11257.7276
$ws.Range("H132").Value = 1283.6538
$ws.Range("I132").Value = 773.4375
$ws.Range("J132").Value = 2100
$ws.Range("K132").Value = 6960.9375
$ws.Range("L132").Value = 18900
$ws.Range("M132").Value = -4430.9375
$ws.Range("N132").Value = -23960
$ws.Range("H135").Value = 1491.0968
$ws.Range("I135").Value = 288.35
$ws.Range("J135").Value = 3677.9092
$ws.Range("K135").Value = 2595.15
$ws.Range("L135").Value = 33101.1828
$ws.Range("M135").Value = -60.15000000000009
$ws.Range("N135").Value = -38171.1828

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H54").Value = 10000
$ws.Range("J54").Value = 10000
$ws.Range("L54").Value = 10000
$ws.Range("N54").Value = -10780
$ws.Range("H118").Value = 14300
$ws.Range("J118").Value = 14300
$ws.Range("L118").Value = 14300
$ws.Range("N118").Value = -17614

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 29129.5
$ws.Range("I29").Value = 28000
$ws.Range("J29").Value = 29506
$ws.Range("K29").Value = 28000
$ws.Range("L29").Value = 29506
$ws.Range("M29").Value = -27705
$ws.Range("N29").Value = -30096
$ws.Range("H32").Value = 9173.571
$ws.Range("I32").Value = 3880
$ws.Range("J32").Value = 22407.5
$ws.Range("K32").Value = 3880
$ws.Range("L32").Value = 22407.5
$ws.Range("M32").Value = -3563
$ws.Range("N32").Value = -23041.5
$ws.Range("H68").Value = 2180.3572
$ws.Range("I68").Value = 2337.4285
$ws.Range("J68").Value = 2023.2858
$ws.Range("K68").Value = 2337.4285
$ws.Range("L68").Value = 2023.2858
$ws.Range("M68").Value = -1588.4285
$ws.Range("N68").Value = -3521.2858
$ws.Range("H71").Value = 2180.3572
$ws.Range("I71").Value = 2337.4285
$ws.Range("J71").Value = 2023.2858
$ws.Range("K71").Value = 11687.1425
$ws.Range("L71").Value = 10116.429
$ws.Range("M71").Value = -7943.1425
$ws.Range("N71").Value = -17604.429

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 21432.223
$ws.Range("I39").Value = 20303.334
$ws.Range("J39").Value = 21996.666
$ws.Range("K39").Value = 20303.334
$ws.Range("L39").Value = 21996.666
$ws.Range("M39").Value = -19890.334
$ws.Range("N39").Value = -22822.666
$ws.Range("H127").Value = 34007.25
$ws.Range("J127").Value = 34007.25
$ws.Range("L127").Value = 34007.25
$ws.Range("N127").Value = -43927.25
$ws.Range("H132").Value = 1126.98
$ws.Range("I132").Value = 1078.3939
$ws.Range("J132").Value = 1221.2941
$ws.Range("K132").Value = 3235.1817
$ws.Range("L132").Value = 3663.8823
$ws.Range("M132").Value = -705.1817000000001
$ws.Range("N132").Value = -8723.882300000001
$ws.Range("H136").Value = 1475.72
$ws.Range("I136").Value = 1338.2222
$ws.Range("J136").Value = 1829.2858
$ws.Range("K136").Value = 4014.6666
$ws.Range("L136").Value = 5487.857400000001
$ws.Range("M136").Value = -1464.6666
$ws.Range("N136").Value = -10587.8574
